$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F:V) between row 54 and row 55 ---
# Row 54 originally: Cagliari vs AC Milan
# Row 55 originally: Empoli vs Salernitana
# After the edit they should be swapped so row 54 becomes Empoli vs Salernitana
# and row 55 becomes Cagliari vs AC Milan (columns A:E stay untouched).

$row54F = $ws.Range("F54").Value2
$row54H = $ws.Range("H54").Value2
$row54I = $ws.Range("I54").Value2
$row54J = $ws.Range("J54").Value2
$row54K = $ws.Range("K54").Value2
$row54L = $ws.Range("L54").Value2
$row54M = $ws.Range("M54").Value2
$row54N = $ws.Range("N54").Value2
$row54O = $ws.Range("O54").Value2
$row54P = $ws.Range("P54").Value2
$row54Q = $ws.Range("Q54").Value2
$row54R = $ws.Range("R54").Value2
$row54S = $ws.Range("S54").Value2
$row54T = $ws.Range("T54").Value2
$row54U = $ws.Range("U54").Value2
$row54V = $ws.Range("V54").Value2

$row55F = $ws.Range("F55").Value2
$row55H = $ws.Range("H55").Value2
$row55I = $ws.Range("I55").Value2
$row55J = $ws.Range("J55").Value2
$row55K = $ws.Range("K55").Value2
$row55L = $ws.Range("L55").Value2
$row55M = $ws.Range("M55").Value2
$row55N = $ws.Range("N55").Value2
$row55O = $ws.Range("O55").Value2
$row55P = $ws.Range("P55").Value2
$row55Q = $ws.Range("Q55").Value2
$row55R = $ws.Range("R55").Value2
$row55S = $ws.Range("S55").Value2
$row55T = $ws.Range("T55").Value2
$row55U = $ws.Range("U55").Value2
$row55V = $ws.Range("V55").Value2

$ws.Range("F54").Value = $row55F
$ws.Range("H54").Value = $row55H
$ws.Range("I54").Value = $row55I
$ws.Range("J54").Value = $row55J
$ws.Range("K54").Value = $row55K
$ws.Range("L54").Value = $row55L
$ws.Range("M54").Value = $row55M
$ws.Range("N54").Value = $row55N
$ws.Range("O54").Value = $row55O
$ws.Range("P54").Value = $row55P
$ws.Range("Q54").Value = $row55Q
$ws.Range("R54").Value = $row55R
$ws.Range("S54").Value = $row55S
$ws.Range("T54").Value = $row55T
$ws.Range("U54").Value = $row55U
$ws.Range("V54").Value = $row55V

$ws.Range("F55").Value = $row54F
$ws.Range("H55").Value = $row54H
$ws.Range("I55").Value = $row54I
$ws.Range("J55").Value = $row54J
$ws.Range("K55").Value = $row54K
$ws.Range("L55").Value = $row54L
$ws.Range("M55").Value = $row54M
$ws.Range("N55").Value = $row54N
$ws.Range("O55").Value = $row54O
$ws.Range("P55").Value = $row54P
$ws.Range("Q55").Value = $row54Q
$ws.Range("R55").Value = $row54R
$ws.Range("S55").Value = $row54S
$ws.Range("T55").Value = $row54T
$ws.Range("U55").Value = $row54U
$ws.Range("V55").Value = $row54V

# --- Append new row 106 with the Verona vs Monza match ---

$ws.Range("A105").Copy()
$ws.Range("A106").PasteSpecial(-4122)
$ws.Range("E105").Copy()
$ws.Range("E106").PasteSpecial(-4122)

$ws.Range("A106").Value = 105
$ws.Range("B106").Value = "italy"
$ws.Range("C106").Value = "serie-a"
$ws.Range("D106").Value = "2023-2024"
$ws.Range("E106").Value = 45235.52083333334
$ws.Range("F106").Value = "Verona"
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = "Monza"
$ws.Range("I106").Value = 3
$ws.Range("J106").Value = 2.58
$ws.Range("K106").Value = "22/10/2023 12:02"
$ws.Range("L106").Value = 3.21
$ws.Range("M106").Value = "05/11/2023 12:17"
$ws.Range("N106").Value = 3.3
$ws.Range("O106").Value = "22/10/2023 12:02"
$ws.Range("P106").Value = 3.15
$ws.Range("Q106").Value = "05/11/2023 12:25"
$ws.Range("R106").Value = 2.74
$ws.Range("S106").Value = "22/10/2023 12:02"
$ws.Range("T106").Value = 2.51
$ws.Range("U106").Value = "05/11/2023 12:26"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/italy/serie-a/verona-monza/2LCmmY2T/"
